# Add a new restaurant ("Baoburg") to the restaurants sheet, fill in a
# missing address for the existing "Bernie's" row, and lightly restyle a
# couple of cells (a muted placeholder-style address font + a phone-number
# hyperlink matching the existing Oxomoco one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("restaurants")

# --- New row 4: Baoburg -----------------------------------------------
# Values are written in the same order the source workbook's shared
# strings were appended in, so the resulting file lines up cell-for-cell.
$ws.Range("I4").Value = " (718) 349-0011"
$ws.Range("C4").Value = "Southern Asian"
$ws.Range("D4").Value = "614 Manhattan Ave, Brooklyn, NY 11222"
$ws.Range("A4").Value = "Baoburg"
$ws.Range("E4").Value = 15
$ws.Range("G4").Value = 4.6
$ws.Range("H4").Value = 20
$ws.Range("J4").Value = "Greenpoint"

# --- Row 3 (Bernie's): fill in the address + table count ---------------
$ws.Range("D3").Value = "332 Driggs Ave, Brooklyn, New York, 11222"
$ws.Range("I3").Value = 69

# --- Formatting ----------------------------------------------------------
# D4 gets a muted, larger placeholder-style font (Arial 14, light grey).
$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 14
$ws.Range("D4").Font.Color = 13025725

# I4 becomes a clickable phone-number lookup, matching I2's existing
# hyperlink style/behaviour.
$url = "https://www.google.com/search?q=baoberg&oq=baoberg&gs_lcrp=EgZjaHJvbWUyBggAEEUYOTIPCAEQLhgKGK8BGMcBGIAEMgkIAhAAGAoYgAQyCQgDEAAYChiABDIJCAQQABgKGIAEMgkIBRAAGAoYgAQyCQgGEAAYChiABDIMCAcQABgFGAoYDxge0gEIMTY5MGowajSoAgCwAgA&sourceid=chrome&ie=UTF-8"
$ws.Hyperlinks.Add($ws.Range("I4"), $url, "", "", $url) | Out-Null
$ws.Range("I4").Style = "Hyperlink"

# New row is a touch taller to accommodate the bigger font.
$ws.Rows.Item(4).RowHeight = 18

# Leave the cursor where the author left it.
$ws.Range("I6").Select() | Out-Null
